$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 38.718679
$ws.Range("H2").Value = 116.156037
$ws.Range("I2").Value = 0.01404461724059496
$ws.Range("J2").Value = 0.01404461724059496
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.93036733333333
$ws.Range("N2").Value = 44.791102
$ws.Range("O2").Value = 0.252612808865421
$ws.Range("P2").Value = 0.252612808865421
$ws.Range("Q2").Value = 578.0841001314193
$ws.Range("R2").Value = 5202.756901182774
$ws.Range("S2").Value = 0.003547850210586412
$ws.Range("T2").Value = 0.003547850210586412

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 38.718679
$ws.Range("H3").Value = 116.156037
$ws.Range("I3").Value = 0.01404461724059496
$ws.Range("J3").Value = 0.01404461724059496
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.28486166666667
$ws.Range("N3").Value = 60.854585
$ws.Range("O3").Value = 0.3432076230048887
$ws.Range("P3").Value = 0.3432076230048887
$ws.Range("Q3").Value = 785.4030474310717
$ws.Range("R3").Value = 7068.627426879645
$ws.Range("S3").Value = 0.004820219699158078
$ws.Range("T3").Value = 0.004820219699158078

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 38.718679
$ws.Range("H4").Value = 116.156037
$ws.Range("I4").Value = 0.01404461724059496
$ws.Range("J4").Value = 0.01404461724059496
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.48898366666667
$ws.Range("N4").Value = 58.466951
$ws.Range("O4").Value = 0.3297418473407271
$ws.Range("P4").Value = 0.3297418473407271
$ws.Range("Q4").Value = 754.5877026259096
$ws.Range("R4").Value = 6791.289323633187
$ws.Range("S4").Value = 0.004631098034107209
$ws.Range("T4").Value = 0.004631098034107209

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 38.718679
$ws.Range("H5").Value = 116.156037
$ws.Range("I5").Value = 0.01404461724059496
$ws.Range("J5").Value = 0.01404461724059496
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.399549333333334
$ws.Range("N5").Value = 13.198648
$ws.Range("O5").Value = 0.07443772078896321
$ws.Range("P5").Value = 0.07443772078896321
$ws.Range("Q5").Value = 170.3447383819974
$ws.Range("R5").Value = 1533.102645437976
$ws.Range("S5").Value = 0.001045449296743267
$ws.Range("T5").Value = 0.001045449296743267

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2511.398112333333
$ws.Range("H6").Value = 7534.194336999999
$ws.Range("I6").Value = 0.9109718135392577
$ws.Range("J6").Value = 0.9109718135392579
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.93036733333333
$ws.Range("N6").Value = 44.791102
$ws.Range("O6").Value = 0.252612808865421
$ws.Range("P6").Value = 0.252612808865421
$ws.Range("Q6").Value = 37496.09633737658
$ws.Range("R6").Value = 337464.8670363893
$ws.Range("S6").Value = 0.2301231486153784
$ws.Range("T6").Value = 0.2301231486153785

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2511.398112333333
$ws.Range("H7").Value = 7534.194336999999
$ws.Range("I7").Value = 0.9109718135392577
$ws.Range("J7").Value = 0.9109718135392579
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.28486166666667
$ws.Range("N7").Value = 60.854585
$ws.Range("O7").Value = 0.3432076230048887
$ws.Range("P7").Value = 0.3432076230048887
$ws.Range("Q7").Value = 50943.36329860945
$ws.Range("R7").Value = 458490.2696874851
$ws.Range("S7").Value = 0.3126524707492613
$ws.Range("T7").Value = 0.3126524707492614

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2511.398112333333
$ws.Range("H8").Value = 7534.194336999999
$ws.Range("I8").Value = 0.9109718135392577
$ws.Range("J8").Value = 0.9109718135392579
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 19.48898366666667
$ws.Range("N8").Value = 58.466951
$ws.Range("O8").Value = 0.3297418473407271
$ws.Range("P8").Value = 0.3297418473407271
$ws.Range("Q8").Value = 48944.59679176182
$ws.Range("R8").Value = 440501.3711258564
$ws.Range("S8").Value = 0.3003855286717672
$ws.Range("T8").Value = 0.3003855286717673

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2511.398112333333
$ws.Range("H9").Value = 7534.194336999999
$ws.Range("I9").Value = 0.9109718135392577
$ws.Range("J9").Value = 0.9109718135392579
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.399549333333334
$ws.Range("N9").Value = 13.198648
$ws.Range("O9").Value = 0.07443772078896321
$ws.Range("P9").Value = 0.07443772078896321
$ws.Range("Q9").Value = 11049.01989085071
$ws.Range("R9").Value = 99441.17901765637
$ws.Range("S9").Value = 0.06781066550285073
$ws.Range("T9").Value = 0.06781066550285074

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.092134333333333
$ws.Range("H10").Value = 15.276403
$ws.Range("I10").Value = 0.001847094980935658
$ws.Range("J10").Value = 0.001847094980935659
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.93036733333333
$ws.Range("N10").Value = 44.791102
$ws.Range("O10").Value = 0.252612808865421
$ws.Range("P10").Value = 0.252612808865421
$ws.Range("Q10").Value = 76.0274361073451
$ws.Range("R10").Value = 684.2469249661059
$ws.Range("S10").Value = 0.0004665998513753779
$ws.Range("T10").Value = 0.0004665998513753779

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 5.092134333333333
$ws.Range("H11").Value = 15.276403
$ws.Range("I11").Value = 0.001847094980935658
$ws.Range("J11").Value = 0.001847094980935659
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 20.28486166666667
$ws.Range("N11").Value = 60.854585
$ws.Range("O11").Value = 0.3432076230048887
$ws.Range("P11").Value = 0.3432076230048887
$ws.Range("Q11").Value = 103.2932405397506
$ws.Range("R11").Value = 929.639164857755
$ws.Range("S11").Value = 0.0006339370778711876
$ws.Range("T11").Value = 0.0006339370778711876

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 5.092134333333333
$ws.Range("H12").Value = 15.276403
$ws.Range("I12").Value = 0.001847094980935658
$ws.Range("J12").Value = 0.001847094980935659
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.48898366666667
$ws.Range("N12").Value = 58.466951
$ws.Range("O12").Value = 0.3297418473407271
$ws.Range("P12").Value = 0.3297418473407271
$ws.Range("Q12").Value = 99.24052285080587
$ws.Range("R12").Value = 893.1647056572531
$ws.Range("S12").Value = 0.0006090645112275091
$ws.Range("T12").Value = 0.0006090645112275091

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 5.092134333333333
$ws.Range("H13").Value = 15.276403
$ws.Range("I13").Value = 0.001847094980935658
$ws.Range("J13").Value = 0.001847094980935659
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.399549333333334
$ws.Range("N13").Value = 13.198648
$ws.Range("O13").Value = 0.07443772078896321
$ws.Range("P13").Value = 0.07443772078896321
$ws.Range("Q13").Value = 22.40309621146045
$ws.Range("R13").Value = 201.627865903144
$ws.Range("S13").Value = 0.0001374935404615839
$ws.Range("T13").Value = 0.0001374935404615839

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 201.6251223333333
$ws.Range("H14").Value = 604.875367
$ws.Range("I14").Value = 0.07313647423921157
$ws.Range("J14").Value = 0.07313647423921157
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 14.93036733333333
$ws.Range("N14").Value = 44.791102
$ws.Range("O14").Value = 0.252612808865421
$ws.Range("P14").Value = 0.252612808865421
$ws.Range("Q14").Value = 3010.337140064937
$ws.Range("R14").Value = 27093.03426058443
$ws.Range("S14").Value = 0.01847521018808074
$ws.Range("T14").Value = 0.01847521018808074

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 201.6251223333333
$ws.Range("H15").Value = 604.875367
$ws.Range("I15").Value = 0.07313647423921157
$ws.Range("J15").Value = 0.07313647423921157
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 20.28486166666667
$ws.Range("N15").Value = 60.854585
$ws.Range("O15").Value = 0.3432076230048887
$ws.Range("P15").Value = 0.3432076230048887
$ws.Range("Q15").Value = 4089.937715056411
$ws.Range("R15").Value = 36809.43943550769
$ws.Range("S15").Value = 0.02510099547859808
$ws.Range("T15").Value = 0.02510099547859808

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 201.6251223333333
$ws.Range("H16").Value = 604.875367
$ws.Range("I16").Value = 0.07313647423921157
$ws.Range("J16").Value = 0.07313647423921157
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 19.48898366666667
$ws.Range("N16").Value = 58.466951
$ws.Range("O16").Value = 0.3297418473407271
$ws.Range("P16").Value = 0.3297418473407271
$ws.Range("Q16").Value = 3929.468715944002
$ws.Range("R16").Value = 35365.21844349601
$ws.Range("S16").Value = 0.02411615612362512
$ws.Range("T16").Value = 0.02411615612362512

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 201.6251223333333
$ws.Range("H17").Value = 604.875367
$ws.Range("I17").Value = 0.07313647423921157
$ws.Range("J17").Value = 0.07313647423921157
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.399549333333334
$ws.Range("N17").Value = 13.198648
$ws.Range("O17").Value = 0.07443772078896321
$ws.Range("P17").Value = 0.07443772078896321
$ws.Range("Q17").Value = 887.0596725448686
$ws.Range("R17").Value = 7983.537052903816
$ws.Range("S17").Value = 0.005444112448907632
$ws.Range("T17").Value = 0.005444112448907632
